# Update column B formulas on the "Sciences" sheet: change the column
# referenced inside SEARCH()/RIGHT() from N to M, for rows 2 through 74.
# Letting Excel recalc afterwards updates the cached <v> values to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sciences")

for ($row = 2; $row -le 74; $row++) {
    $formula = '=IF(OR(NOT(ISERROR(SEARCH("archive.org",M' + $row + '))),NOT(ISERROR(SEARCH("app.box.com",M' + $row + '))),NOT(ISERROR(SEARCH("islamway.net",M' + $row + '))),NOT(ISERROR(SEARCH("qurancomplex.gov.sa",M' + $row + '))),NOT(ISERROR(SEARCH("tanzil.net",M' + $row + '))),NOT(ISERROR(SEARCH("alsirah.com",M' + $row + '))),NOT(ISERROR(SEARCH("i36",M' + $row + '))),(RIGHT(M' + $row + ',4)=".pdf"),C' + $row + '=6,C' + $row + '=8,C' + $row + '=9),0,1)'
    $ws.Range("B$row").Formula = $formula
}

$excel.CalculateFullRebuild()
$wb.Save()
